$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Alemania
$ws.Range("B8").Value = 138135
$ws.Range("C8").Value = 437
$ws.Range("E8").Value = 52242
$ws.Range("G8").Value = 41
$ws.Range("H8").Value = 4093

# Row 20 - Austria
$ws.Range("B20").Value = 14508
$ws.Range("C20").Value = 32
$ws.Range("E20").Value = 5112

# Row 30 - Polonia
$ws.Range("B30").Value = 8214
$ws.Range("C30").Value = 296
$ws.Range("E30").Value = 7030
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 318

# Row 32 - was Pakistan, becomes Dinamarca with new totals
$ws.Range("A32").Value = "Dinamarca"
$ws.Range("B32").Value = 7073
$ws.Range("C32").Value = 194
$ws.Range("D32").Value = 3023
$ws.Range("E32").Value = 3729
$ws.Range("F32").Value = 92
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 321

# Row 33 - was Noruega, becomes Pakistan (old Pakistan values shift down)
$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 7025
$ws.Range("C33").Value = 106
$ws.Range("D33").Value = 1765
$ws.Range("E33").Value = 5125
$ws.Range("F33").Value = 46
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 135

# Row 34 - was Dinamarca, becomes Noruega (old Noruega values shift down)
$ws.Range("A34").Value = "Noruega"
$ws.Range("B34").Value = 6905
$ws.Range("D34").Value = 32
$ws.Range("E34").Value = 6721
$ws.Range("F34").Value = 64
$ws.Range("H34").Value = 152

# Row 39 - was Emiratos Arabes Unidos, becomes Filipinas with new totals
$ws.Range("A39").Value = "Filipinas"
$ws.Range("B39").Value = 5878
$ws.Range("C39").Value = 218
$ws.Range("D39").Value = 487
$ws.Range("E39").Value = 5004
$ws.Range("G39").Value = 25
$ws.Range("H39").Value = 387

# Row 40 - was Filipinas, becomes Emiratos Arabes Unidos (old Emiratos values shift down)
$ws.Range("A40").Value = "Emiratos Arabes Unidos"
$ws.Range("B40").Value = 5825
$ws.Range("D40").Value = 1095
$ws.Range("E40").Value = 4695
$ws.Range("H40").Value = 35

# Row 71 - Uzbekistan
$ws.Range("D71").Value = 140
$ws.Range("E71").Value = 1236

# Row 78 - Oman
$ws.Range("B78").Value = 1069
$ws.Range("C78").Value = 50
$ws.Range("E78").Value = 888

# Row 89 - Letonia
$ws.Range("D89").Value = 88
$ws.Range("E89").Value = 589
